$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table with the latest scrape.
# Only cells whose content actually changed are touched; for numeric-
# looking Price strings we force text format first so Excel keeps them
# as literal text (e.g. "1.00") instead of silently converting to a number.

$ws.Range("D2").Value = "63.093.99"
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").Value = "2.681.88"
$ws.Range("E3").Value = "  -2.19%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.79"
$ws.Range("E5").Value = "  -2.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.93"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("E9").Value = "  -2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("E11").Value = "  -3.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("E12").Value = "  -6.68%  "

$ws.Range("D13").Value = "3.153.94"
$ws.Range("E13").Value = "  -2.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.29"
$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").Value = "62.862.41"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("D17").Value = "2.679.50"
$ws.Range("E17").Value = "  -2.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.90"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.62"
$ws.Range("E19").Value = "  -3.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.86"
$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("E21").Value = "  -3.68%  "

$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.508"
$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.18"
$ws.Range("E24").Value = "  -1.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.18"
$ws.Range("E27").Value = "  -2.63%  "

$ws.Range("D28").Value = "0.0₃0866"
$ws.Range("E28").Value = "  -4.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.38"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.27"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.96"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.60"
$ws.Range("E32").Value = "  +2.24%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.49"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.57"
$ws.Range("E36").Value = "  -2.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "350.22"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.40"
$ws.Range("E39").Value = "  +2.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.969"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.98"
$ws.Range("E41").Value = "  -3.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.18"
$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.46"
$ws.Range("E43").Value = "  -3.09%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.88"
$ws.Range("E44").Value = "  -5.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0566"
$ws.Range("E45").Value = "  -2.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.619"
$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.05"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0975"
$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0242"
$ws.Range("E50").Value = "  -1.95%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.27"
$ws.Range("E51").Value = "  -3.96%  "
